$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: A1 "File Name" -> "Loc", B1 "Unnormalized P_max" -> "P_max"
$ws.Range("A1").Value2 = "Loc"
$ws.Range("B1").Value2 = "P_max"

# Replace column A values (filenames) with the corresponding electrode
# location codes that were previously stored in column C, for data rows 2-67.
for ($r = 2; $r -le 67; $r++) {
    $loc = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r, 1).Value2 = $loc
}

# Remove the now-redundant "Electrode Locations" column C entirely.
$ws.Columns.Item(3).Delete()
